# Auto-generated edit script applying value changes per diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.76
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 4.8
$ws.Range("K2").Value = 4.8
$ws.Range("Q2").Value = 1.54
$ws.Range("V2").Value = 1.26
$ws.Range("Z2").Value = 1000
$ws.Range("AB2").Value = 13.5
$ws.Range("AD2").Value = 18.5
$ws.Range("AE2").Value = 50
$ws.Range("AI2").Value = 1000
$ws.Range("AL2").Value = 27
$ws.Range("AM2").Value = 1000
# Row 3
$ws.Range("F3").Value = 1.45
$ws.Range("G3").Value = 1.47
$ws.Range("H3").Value = 8
$ws.Range("I3").Value = 8.6
$ws.Range("J3").Value = 5.1
$ws.Range("K3").Value = 5.3
$ws.Range("P3").Value = 2.48
$ws.Range("Q3").Value = 1.64
$ws.Range("R3").Value = 1.59
$ws.Range("S3").Value = 2.6
$ws.Range("T3").Value = 1.87
$ws.Range("U3").Value = 2.08
$ws.Range("V3").Value = 1.13
$ws.Range("W3").Value = 3.1
$ws.Range("X3").Value = 21
$ws.Range("Y3").Value = 32
$ws.Range("Z3").Value = 75
$ws.Range("AA3").Value = 260
$ws.Range("AC3").Value = 11.5
$ws.Range("AD3").Value = 29
$ws.Range("AE3").Value = 110
$ws.Range("AF3").Value = 9
$ws.Range("AG3").Value = 9.800000000000001
$ws.Range("AH3").Value = 23
$ws.Range("AI3").Value = 95
$ws.Range("AJ3").Value = 12.5
$ws.Range("AN3").Value = 5.7
$ws.Range("AO3").Value = 110
# Row 4
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 5.1
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 2.42
$ws.Range("Q4").Value = 1.66
$ws.Range("R4").Value = 1.56
$ws.Range("S4").Value = 2.66
$ws.Range("T4").Value = 1.63
$ws.Range("U4").Value = 2.52
$ws.Range("X4").Value = 21
$ws.Range("AA4").Value = 80
$ws.Range("AB4").Value = 12.5
$ws.Range("AC4").Value = 9.199999999999999
$ws.Range("AF4").Value = 14
# Row 5
$ws.Range("F5").Value = 1.9
$ws.Range("G5").Value = 3.4
$ws.Range("H5").Value = 1.71
$ws.Range("I5").Value = 2.94
$ws.Range("J5").Value = 3.6
$ws.Range("K5").Value = 950
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 2.2
$ws.Range("O5").Value = 1.18
$ws.Range("P5").Value = 2.18
$ws.Range("R5").Value = 1.49
$ws.Range("S5").Value = 2.22
$ws.Range("T5").Value = 1.04
$ws.Range("U5").Value = 1.04
$ws.Range("V5").Value = 1.51
$ws.Range("W5").Value = 1.42
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000
# Row 6
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 2.12
$ws.Range("O6").Value = 1.25
$ws.Range("R6").Value = 1.13
$ws.Range("S6").Value = 1.71
$ws.Range("T6").Value = 1.01
$ws.Range("U6").Value = 1.01
$ws.Range("V6").Value = 1.12
$ws.Range("W6").Value = 2.6
$ws.Range("X6").Value = 29
$ws.Range("Y6").Value = 38
$ws.Range("Z6").Value = 90
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 13
$ws.Range("AC6").Value = 15.5
$ws.Range("AD6").Value = 40
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 14.5
$ws.Range("AG6").Value = 15
$ws.Range("AH6").Value = 32
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 20
$ws.Range("AK6").Value = 23
$ws.Range("AL6").Value = 50
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 10.5
$ws.Range("AO6").Value = 1000
# Row 7
$ws.Range("F7").Value = 2.44
$ws.Range("G7").Value = 2.48
$ws.Range("H7").Value = 3.35
$ws.Range("K7").Value = 3.4
$ws.Range("L7").Value = 1.44
$ws.Range("N7").Value = 3.5
$ws.Range("Q7").Value = 2.16
$ws.Range("V7").Value = 1.42
$ws.Range("W7").Value = 1.67
$ws.Range("X7").Value = 12
$ws.Range("Z7").Value = 21
$ws.Range("AA7").Value = 60
$ws.Range("AB7").Value = 9.6
$ws.Range("AC7").Value = 7.2
$ws.Range("AD7").Value = 14
$ws.Range("AE7").Value = 40
$ws.Range("AG7").Value = 11.5
$ws.Range("AJ7").Value = 34
$ws.Range("AK7").Value = 28
$ws.Range("AL7").Value = 46
$ws.Range("AM7").Value = 120
$ws.Range("AO7").Value = 42
# Row 8
$ws.Range("G8").Value = 2.64
$ws.Range("H8").Value = 2.64
$ws.Range("I8").Value = 2.66
$ws.Range("L8").Value = 1.25
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 6.2
$ws.Range("O8").Value = 1.17
$ws.Range("V8").Value = 1.6
$ws.Range("W8").Value = 1.6
$ws.Range("X8").Value = 25
$ws.Range("Y8").Value = 18.5
$ws.Range("Z8").Value = 22
$ws.Range("AA8").Value = 40
$ws.Range("AB8").Value = 18.5
$ws.Range("AC8").Value = 9.6
$ws.Range("AD8").Value = 12
$ws.Range("AF8").Value = 22
$ws.Range("AH8").Value = 13.5
$ws.Range("AI8").Value = 26
$ws.Range("AJ8").Value = 40
$ws.Range("AK8").Value = 23
$ws.Range("AL8").Value = 27
$ws.Range("AM8").Value = 48
$ws.Range("AN8").Value = 13
$ws.Range("AO8").Value = 13
# Row 9
$ws.Range("G9").Value = 10.5
$ws.Range("H9").Value = 1.37
$ws.Range("J9").Value = 5.7
$ws.Range("K9").Value = 5.9
$ws.Range("L9").Value = 1.28
$ws.Range("R9").Value = 1.61
$ws.Range("T9").Value = 1.96
$ws.Range("V9").Value = 3.55
$ws.Range("W9").Value = 1.1
$ws.Range("X9").Value = 23
$ws.Range("Z9").Value = 8.6
$ws.Range("AD9").Value = 9.800000000000001
$ws.Range("AE9").Value = 13.5
$ws.Range("AG9").Value = 36
$ws.Range("AH9").Value = 26
$ws.Range("AI9").Value = 32
$ws.Range("AJ9").Value = 360
$ws.Range("AK9").Value = 150
$ws.Range("AM9").Value = 140
$ws.Range("AN9").Value = 150
# Row 10
$ws.Range("H10").Value = 19.5
$ws.Range("L10").Value = 1.25
$ws.Range("O10").Value = 1.16
$ws.Range("P10").Value = 2.92
$ws.Range("Q10").Value = 1.49
$ws.Range("R10").Value = 1.77
$ws.Range("S10").Value = 2.26
$ws.Range("T10").Value = 2.32
$ws.Range("V10").Value = 1.04
$ws.Range("W10").Value = 6
$ws.Range("Y10").Value = 65
$ws.Range("AC10").Value = 20
$ws.Range("AF10").Value = 7.8
$ws.Range("AH10").Value = 50
$ws.Range("AL10").Value = 160
$ws.Range("AN10").Value = 3.35
# Row 11
$ws.Range("H11").Value = 11
$ws.Range("J11").Value = 7
$ws.Range("K11").Value = 7.2
$ws.Range("L11").Value = 1.19
$ws.Range("P11").Value = 3.45
$ws.Range("S11").Value = 1.94
$ws.Range("V11").Value = 1.09
$ws.Range("W11").Value = 4.2
$ws.Range("X11").Value = 42
$ws.Range("Z11").Value = 130
$ws.Range("AA11").Value = 380
$ws.Range("AB11").Value = 15
$ws.Range("AC11").Value = 16
$ws.Range("AD11").Value = 40
$ws.Range("AE11").Value = 140
$ws.Range("AF11").Value = 10.5
$ws.Range("AG11").Value = 11
$ws.Range("AH11").Value = 25
$ws.Range("AJ11").Value = 12
$ws.Range("AK11").Value = 12
$ws.Range("AL11").Value = 26
$ws.Range("AM11").Value = 100
$ws.Range("AN11").Value = 3.45
$ws.Range("AO11").Value = 110
# Row 12
$ws.Range("F12").Value = 6.2
$ws.Range("G12").Value = 6.4
$ws.Range("H12").Value = 1.63
$ws.Range("I12").Value = 1.65
$ws.Range("J12").Value = 4.3
$ws.Range("L12").Value = 1.36
$ws.Range("U12").Value = 2.04
$ws.Range("V12").Value = 2.54
$ws.Range("W12").Value = 1.18
$ws.Range("X12").Value = 17
$ws.Range("Z12").Value = 9.6
$ws.Range("AC12").Value = 9.4
$ws.Range("AD12").Value = 9.6
$ws.Range("AE12").Value = 16
$ws.Range("AF12").Value = 50
$ws.Range("AG12").Value = 23
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 34
$ws.Range("AJ12").Value = 170
$ws.Range("AK12").Value = 85
$ws.Range("AM12").Value = 120
$ws.Range("AN12").Value = 100
$ws.Range("AO12").Value = 8.4
# Row 13
$ws.Range("F13").Value = 3.05
$ws.Range("H13").Value = 2.42
$ws.Range("I13").Value = 2.46
$ws.Range("L13").Value = 1.33
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 4.8
$ws.Range("V13").Value = 1.68
$ws.Range("W13").Value = 1.46
$ws.Range("X13").Value = 17.5
$ws.Range("Z13").Value = 17
$ws.Range("AC13").Value = 8.4
$ws.Range("AD13").Value = 11
$ws.Range("AE13").Value = 23
$ws.Range("AG13").Value = 13
$ws.Range("AI13").Value = 32
$ws.Range("AJ13").Value = 50
$ws.Range("AN13").Value = 23

Write-Output "Applied 265 cell updates"
